$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77-126 down to 78-127
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly Espárragos record
$ws.Cells.Item(77, 1).Value = 5
$ws.Cells.Item(77, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(77, 3).Value = "Maule"
$ws.Cells.Item(77, 4).Value = 45233
$ws.Cells.Item(77, 5).Value = 7
$ws.Cells.Item(77, 6).Value = 300000000
$ws.Cells.Item(77, 7).Value = "Espárragos"
$ws.Cells.Item(77, 8).Value = "Verde"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 3000
$ws.Cells.Item(77, 11).Value = 1300
$ws.Cells.Item(77, 12).Value = 1300
$ws.Cells.Item(77, 13).Value = 1300
$ws.Cells.Item(77, 14).Value = "`$/kilo"
$ws.Cells.Item(77, 15).Value = "Provincia de Linares"
$ws.Cells.Item(77, 16).Value = 1300
$ws.Cells.Item(77, 17).Value = 1
$ws.Cells.Item(77, 18).Value = "Hortaliza"
